# Auto-generated Word COM-interop script.
# Applies 5 paragraph-level OOXML replacements via Range.InsertXML,
# matching the target diff (run-splitting for proofErr spellcheck/grammar
# marks, plus an inserted oMath accent formula).

$d = $word.ActiveDocument

# --- edit 1: find unique text, grab its paragraph, splice replacement XML ---
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Tb ch 3.2", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) { throw "edit 1: Find failed for Tb ch 3.2" }
$para1 = $rng1.Paragraphs(1).Range
$xml1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p w14:paraId="77DCA286" w14:textId="7F731E7F" w:rsidR="00F61173" w:rsidRDefault="009064E4" w:rsidP="00F61173"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="009064E4"><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Linear Filters</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="000470C3"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r><w:r w:rsidR="000470C3" w:rsidRPr="000470C3"><w:rPr><w:highlight w:val="magenta"/><w:lang w:val="en-US"/></w:rPr><w:t>TODO</w:t></w:r><w:r w:rsidR="000470C3"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r w:rsidR="00F61173"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Tb </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>ch</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> 3.2</w:t></w:r><w:r w:rsidR="000470C3"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r></w:p></w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@
$para1.InsertXML($xml1) | Out-Null

# --- edit 2: find unique text, grab its paragraph, splice replacement XML ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Correlation (cv2.filter2D", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "edit 2: Find failed for Correlation (cv2.filter2D" }
$para2 = $rng2.Paragraphs(1).Range
$xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p w14:paraId="76205AE2" w14:textId="0DA6517B" w:rsidR="00E42A2A" w:rsidRDefault="00915119" w:rsidP="00915119"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="006341C5"><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:highlight w:val="yellow"/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>Correlation</w:t></w:r><w:r w:rsidR="00D40AA2"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> (</w:t></w:r><w:r w:rsidR="007C229B"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">cv2.filter2D, </w:t></w:r><w:r w:rsidR="00D40AA2"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>2D moving average</w:t></w:r><w:r w:rsidR="00BE07C1"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> with (non-)uniform weights</w:t></w:r><w:r w:rsidR="00D40AA2"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r><w:r w:rsidR="00001D91"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>:</w:t></w:r><w:r w:rsidR="00C20A49"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00D94D22"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>G</w:t></w:r><w:r w:rsidR="00085351"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">iven input </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:lang w:val="en-US"/></w:rPr><m:t>I</m:t></m:r></m:oMath><w:r w:rsidR="00B40CF4"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>,</w:t></w:r><w:r w:rsidR="0063336A"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:lang w:val="en-US"/></w:rPr><m:t>G=F⊗I</m:t></m:r></m:oMath><w:r w:rsidR="0063336A"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>where</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p></w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@
$para2.InsertXML($xml2) | Out-Null

# --- edit 3: find unique text, grab its paragraph, splice replacement XML ---
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("exact match of image crop and filter results in 1.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) { throw "edit 3: Find failed for exact match of image crop and filter results in 1." }
$para3 = $rng3.Paragraphs(1).Range
$xml3 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p w14:paraId="20CDBE2F" w14:textId="4D7363DF" w:rsidR="000F004A" w:rsidRDefault="00004A20" w:rsidP="00915119"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00004A20"><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t>Normalized Cross-Correlation</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>:</w:t></w:r><w:r w:rsidR="0087457D"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> exact match of image crop and filter results in 1.</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Normalized prevents </w:t></w:r><m:oMath><m:acc><m:accPr><m:chr m:val="⃗"/><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/><w:lang w:val="en-US"/></w:rPr></m:ctrlPr></m:accPr><m:e><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/><w:lang w:val="en-US"/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:lang w:val="en-US"/></w:rPr><m:t>t</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:lang w:val="en-US"/></w:rPr><m:t>ij</m:t></m:r></m:sub></m:sSub></m:e></m:acc></m:oMath><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> that is all or almost all white (255) </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>to generate large response.</w:t></w:r></w:p></w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@
$para3.InsertXML($xml3) | Out-Null

# --- edit 4: find unique text, grab its paragraph, splice replacement XML ---
$rng4 = $d.Content
$found4 = $rng4.Find.Execute("Convolution: operator that flips filter horizontal", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found4) { throw "edit 4: Find failed for Convolution: operator that flips filter horizontal" }
$para4 = $rng4.Paragraphs(1).Range
$xml4 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p w14:paraId="23630574" w14:textId="495C3B03" w:rsidR="00915119" w:rsidRDefault="00915119" w:rsidP="00915119"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="006145CC"><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:highlight w:val="yellow"/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>Convolution</w:t></w:r><w:r w:rsidR="006145CC"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r w:rsidR="008B4EC0"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>operator that flips filter horizontal</w:t></w:r><w:r w:rsidR="00420D92"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>ly</w:t></w:r><w:r w:rsidR="008B4EC0"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> and vertical</w:t></w:r><w:r w:rsidR="00420D92"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>ly then applies correlation.</w:t></w:r><w:r w:rsidR="00D11B1F" w:rsidRPr="00D11B1F"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00D11B1F"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Given input </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:lang w:val="en-US"/></w:rPr><m:t>I</m:t></m:r></m:oMath><w:r w:rsidR="00D11B1F"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:lang w:val="en-US"/></w:rPr><m:t>G=F*I</m:t></m:r></m:oMath><w:r w:rsidR="004A5DE4"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>where</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p></w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@
$para4.InsertXML($xml4) | Out-Null

# --- edit 5: find unique text, grab its paragraph, splice replacement XML ---
$rng5 = $d.Content
$found5 = $rng5.Find.Execute("and mult in ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found5) { throw "edit 5: Find failed for and mult in " }
$para5 = $rng5.Paragraphs(1).Range
$xml5 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p w14:paraId="43774FCA" w14:textId="065BEFD3" w:rsidR="000B2D6E" w:rsidRDefault="000044C8"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Method 2: </w:t></w:r><w:r w:rsidR="00510B67"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>FFT</w:t></w:r><w:r w:rsidR="00087297"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> and IFFT run in</w:t></w:r><w:r w:rsidR="00510B67"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:lang w:val="en-US"/></w:rPr><m:t>N</m:t></m:r><m:func><m:funcPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/><w:lang w:val="en-US"/></w:rPr></m:ctrlPr></m:funcPr><m:fName><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:lang w:val="en-US"/></w:rPr><m:t>log</m:t></m:r></m:fName><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:lang w:val="en-US"/></w:rPr><m:t>N</m:t></m:r></m:e></m:func></m:oMath><w:r w:rsidR="00E2677F"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>mult</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> in </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:lang w:val="en-US"/></w:rPr><m:t>N</m:t></m:r></m:oMath><w:r w:rsidR="00B25763"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@
$para5.InsertXML($xml5) | Out-Null

Write-Output "done"
